$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and 1h volume-change (E) figures for the cryptos table.
$updates = @(
    @{ Cell = "D2"; Value = '41.880.88'; ForceText = $true }
    @{ Cell = "E2"; Value = '  -0.93%  '; ForceText = $false }
    @{ Cell = "D3"; Value = '2.212.83'; ForceText = $true }
    @{ Cell = "E3"; Value = '  -1.61%  '; ForceText = $false }
    @{ Cell = "E4"; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = "D5"; Value = '241.42'; ForceText = $true }
    @{ Cell = "E5"; Value = '  -2.22%  '; ForceText = $false }
    @{ Cell = "D6"; Value = '0.626'; ForceText = $true }
    @{ Cell = "E6"; Value = '  +0.54%  '; ForceText = $false }
    @{ Cell = "D7"; Value = '72.77'; ForceText = $true }
    @{ Cell = "E7"; Value = '  -2.63%  '; ForceText = $false }
    @{ Cell = "E8"; Value = '  +0.18%  '; ForceText = $false }
    @{ Cell = "D9"; Value = '0.604'; ForceText = $true }
    @{ Cell = "E9"; Value = '  -2.59%  '; ForceText = $false }
    @{ Cell = "D10"; Value = '42.10'; ForceText = $true }
    @{ Cell = "E10"; Value = '  -1.48%  '; ForceText = $false }
    @{ Cell = "D11"; Value = '0.0953'; ForceText = $true }
    @{ Cell = "E11"; Value = '  +1.15%  '; ForceText = $false }
    @{ Cell = "D12"; Value = '7.03'; ForceText = $true }
    @{ Cell = "E12"; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = "D13"; Value = '0.103'; ForceText = $true }
    @{ Cell = "E13"; Value = '  +0.25%  '; ForceText = $false }
    @{ Cell = "D14"; Value = '2.545.40'; ForceText = $true }
    @{ Cell = "E14"; Value = '  -1.51%  '; ForceText = $false }
    @{ Cell = "D15"; Value = '14.24'; ForceText = $true }
    @{ Cell = "E15"; Value = '  -1.73%  '; ForceText = $false }
    @{ Cell = "D16"; Value = '0.835'; ForceText = $true }
    @{ Cell = "E16"; Value = '  -2.22%  '; ForceText = $false }
    @{ Cell = "D17"; Value = '2.216.32'; ForceText = $true }
    @{ Cell = "E17"; Value = '  -0.60%  '; ForceText = $false }
    @{ Cell = "D18"; Value = '41.841.16'; ForceText = $true }
    @{ Cell = "E18"; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = "D19"; Value = '0.0000106'; ForceText = $true }
    @{ Cell = "E19"; Value = '  +5.59%  '; ForceText = $false }
    @{ Cell = "D20"; Value = '6.18'; ForceText = $true }
    @{ Cell = "E20"; Value = '  +0.60%  '; ForceText = $false }
    @{ Cell = "D21"; Value = '72.71'; ForceText = $true }
    @{ Cell = "E21"; Value = '  +0.40%  '; ForceText = $false }
    @{ Cell = "D22"; Value = '10.62'; ForceText = $true }
    @{ Cell = "E22"; Value = '  +18.58%  '; ForceText = $false }
    @{ Cell = "D23"; Value = '229.95'; ForceText = $true }
    @{ Cell = "E23"; Value = '  -0.79%  '; ForceText = $false }
    @{ Cell = "D24"; Value = '2.08'; ForceText = $true }
    @{ Cell = "E24"; Value = '  -6.46%  '; ForceText = $false }
    @{ Cell = "D25"; Value = '11.83'; ForceText = $true }
    @{ Cell = "E25"; Value = '  +2.96%  '; ForceText = $false }
    @{ Cell = "E26"; Value = '  +0.14%  '; ForceText = $false }
    @{ Cell = "D27"; Value = '3.65'; ForceText = $true }
    @{ Cell = "E27"; Value = '  +0.88%  '; ForceText = $false }
    @{ Cell = "D28"; Value = '2.27'; ForceText = $true }
    @{ Cell = "E28"; Value = '  -1.83%  '; ForceText = $false }
    @{ Cell = "E29"; Value = '  -0.19%  '; ForceText = $false }
    @{ Cell = "D30"; Value = '168.06'; ForceText = $true }
    @{ Cell = "E30"; Value = '  -0.41%  '; ForceText = $false }
    @{ Cell = "D31"; Value = '20.46'; ForceText = $true }
    @{ Cell = "E31"; Value = '  -1.15%  '; ForceText = $false }
    @{ Cell = "D32"; Value = '5.61'; ForceText = $true }
    @{ Cell = "E32"; Value = '  +7.72%  '; ForceText = $false }
    @{ Cell = "D33"; Value = '0.0795'; ForceText = $true }
    @{ Cell = "E33"; Value = '  -3.11%  '; ForceText = $false }
    @{ Cell = "D34"; Value = '29.67'; ForceText = $true }
    @{ Cell = "E34"; Value = '  -3.55%  '; ForceText = $false }
    @{ Cell = "E35"; Value = '  -0.28%  '; ForceText = $false }
    @{ Cell = "D36"; Value = '0.109'; ForceText = $true }
    @{ Cell = "E36"; Value = '  -10.62%  '; ForceText = $false }
    @{ Cell = "D37"; Value = '4.26'; ForceText = $true }
    @{ Cell = "E38"; Value = '  -4.46%  '; ForceText = $false }
    @{ Cell = "D39"; Value = '13.94'; ForceText = $true }
    @{ Cell = "E39"; Value = '  +1.22%  '; ForceText = $false }
    @{ Cell = "D40"; Value = '65.51'; ForceText = $true }
    @{ Cell = "E40"; Value = '  +4.42%  '; ForceText = $false }
    @{ Cell = "D41"; Value = '2.12'; ForceText = $true }
    @{ Cell = "E41"; Value = '  -2.68%  '; ForceText = $false }
    @{ Cell = "D42"; Value = '5.64'; ForceText = $true }
    @{ Cell = "E42"; Value = '  -2.83%  '; ForceText = $false }
    @{ Cell = "E43"; Value = '  -3.13%  '; ForceText = $false }
    @{ Cell = "D44"; Value = '8.79'; ForceText = $true }
    @{ Cell = "E44"; Value = '  +0.96%  '; ForceText = $false }
    @{ Cell = "D45"; Value = '105.19'; ForceText = $true }
    @{ Cell = "E45"; Value = '  -2.41%  '; ForceText = $false }
    @{ Cell = "E46"; Value = '  -2.17%  '; ForceText = $false }
    @{ Cell = "D47"; Value = '2.43'; ForceText = $true }
    @{ Cell = "E47"; Value = '  +5.67%  '; ForceText = $false }
    @{ Cell = "D48"; Value = '1.11'; ForceText = $true }
    @{ Cell = "E48"; Value = '  -0.65%  '; ForceText = $false }
    @{ Cell = "E49"; Value = '  -0.77%  '; ForceText = $false }
    @{ Cell = "E50"; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = "D51"; Value = '2.423.16'; ForceText = $true }
    @{ Cell = "E51"; Value = '  -1.54%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Prefix with a literal apostrophe so numeric-looking strings
        # (e.g. '42.10') are stored as text, not coerced to a Double
        # (which would silently drop the trailing zero / precision).
        $rng.Value = "'" + $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

Write-Host "Updated $($updates.Count) cells on $($ws.Name) with refreshed cryptos data"
